# Checklist name bug fixed:
# Rename the "Completion time" column header to "End time" on the Form1
# worksheet / Table1, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Form1")

# Update the header cell text (this also drives the ListObject/table
# column name, which Excel keeps in sync automatically).
$ws.Range("C1").Value = "End time"

# Restore the saved selection state on the sheet.
$ws.Activate()
$ws.Range("C12").Select()
